$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title strings in row 1 (A1:C1) from "7.2.1" to "7.2.1.1"
$ws.Range("A1").Value = " 7.2.1.1 Энергия керектөөлөрүнүн жалпы көлөмүндөгү энергиянын жаңыланган булактарынын  үлүшү"
$ws.Range("B1").Value = " 7.2.1.1 Доля возобновляемых источников энергии в общем объеме энергопотребления"
$ws.Range("C1").Value = "7.2.1.1 Renewable energy share in the total energy consumption"

# Fix data values
$ws.Range("Q5").Value = 36.700000000000003
$ws.Range("P6").Value = 13859.2
$ws.Range("Q6").Value = 13979.2

# Update the selected cell in the sheet view
$ws.Range("P7").Select()
